# Apply odds updates to Sheet1 as described in the commit's diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("J3").Value = 1.1
$ws.Range("K3").Value = 7

# Row 10
$ws.Range("G10").Value = 1.95
$ws.Range("H10").Value = 3.1
$ws.Range("J10").Value = 1.11
$ws.Range("K10").Value = 6.5
$ws.Range("W10").Value = 17
$ws.Range("AE10").Value = 9

# Row 12
$ws.Range("G12").Value = 1.4
$ws.Range("H12").Value = 4.33
$ws.Range("J12").Value = 1.03
$ws.Range("K12").Value = 15
$ws.Range("L12").Value = 1.2
$ws.Range("M12").Value = 4.33
$ws.Range("N12").Value = 1.67
$ws.Range("O12").Value = 2.15
$ws.Range("R12").Value = 1.91
$ws.Range("S12").Value = 1.8
$ws.Range("W12").Value = 9.5
$ws.Range("AA12").Value = 9
$ws.Range("AE12").Value = 17
$ws.Range("AF12").Value = 41

# Row 13
$ws.Range("G13").Value = 3.25
$ws.Range("H13").Value = 3.1
$ws.Range("I13").Value = 2.15
$ws.Range("J13").Value = 1.08
$ws.Range("K13").Value = 7.5
$ws.Range("N13").Value = 2.3
$ws.Range("O13").Value = 1.6
$ws.Range("T13").Value = 8.5
$ws.Range("U13").Value = 15
$ws.Range("X13").Value = 29
$ws.Range("AA13").Value = 6
$ws.Range("AB13").Value = 17
$ws.Range("AD13").Value = 900
$ws.Range("AE13").Value = 6.5
$ws.Range("AF13").Value = 9.5
$ws.Range("AH13").Value = 21
$ws.Range("AI13").Value = 21

# Row 18
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 3.2
$ws.Range("I18").Value = 3.9
$ws.Range("J18").Value = 1.06
$ws.Range("K18").Value = 10
$ws.Range("L18").Value = 1.33
$ws.Range("M18").Value = 3.25
$ws.Range("N18").Value = 2.08
$ws.Range("O18").Value = 1.73
$ws.Range("R18").Value = 1.91
$ws.Range("S18").Value = 1.91
$ws.Range("U18").Value = 9.5
$ws.Range("W18").Value = 19
$ws.Range("Z18").Value = 8.5
$ws.Range("AA18").Value = 6
$ws.Range("AD18").Value = 301
$ws.Range("AE18").Value = 10
$ws.Range("AI18").Value = 29

# Row 19
$ws.Range("G19").Value = 1.48
$ws.Range("I19").Value = 6.5
$ws.Range("N19").Value = 1.95
$ws.Range("O19").Value = 1.9
$ws.Range("R19").Value = 2.05
$ws.Range("S19").Value = 1.7
$ws.Range("T19").Value = 6
$ws.Range("Z19").Value = 10
$ws.Range("AA19").Value = 8.5
$ws.Range("AD19").Value = 501
$ws.Range("AE19").Value = 15
$ws.Range("AF19").Value = 34
